# Update countries & provincias Spain
# - Refresh COVID data for a handful of countries (Rusia, Filipinas, Singapur,
#   Armenia, Afganistan, Australia, Croacia, Georgia, Estonia, Letonia).
# - Re-sort the whole "Pais" table (rows 4..219) by "Casos totales" (column B)
#   descending, since the sheet is expected to stay ordered that way.
# - Fix up one tie (Santa Lucia / Timor Oriental, both at 27 total cases) whose
#   relative order flips in the source data refresh even though the values
#   are identical.
# - Bump the "Datos actualizados" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 4
$lastDataRow = 219
$nameCol = 1
$firstValCol = 2
$lastValCol = 8

# New per-country figures: Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes.
$updates = @{
    "Rusia"      = @(1159573, 8135, 945920, 193268, 0, 61, 20385)
    "Filipinas"  = @(307288, 3073, 252665, 49242, 0, 37, 5381)
    "Singapur"   = @(57715, 15, 57367, 321, 0, 0, 27)
    "Armenia"    = @(49574, 174, 43665, 4952, 0, 6, 957)
    "Afganistan" = @(39233, 6, 32642, 5136, 0, 2, 1455)
    "Australia"  = @(27044, 4, 24676, 1493, 0, 3, 875)
    "Croacia"    = @(16245, 48, 14793, 1180, 0, 0, 272)
    "Georgia"    = @(5552, 298, 2054, 3467, 0, 3, 31)
    "Estonia"    = @(3267, 67, 2513, 690, 0, 0, 64)
    "Letonia"    = @(1697, 21, 1304, 357, 0, 0, 36)
}

# --- Step 1: apply the data refresh, matching rows by country name ---
foreach ($countryName in $updates.Keys) {
    $newVals = $updates[$countryName]
    for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
        $existingName = $ws.Cells.Item($r, $nameCol).Value2
        if ($existingName -eq $countryName) {
            for ($c = $firstValCol; $c -le $lastValCol; $c++) {
                $ws.Cells.Item($r, $c).Value = $newVals[$c - $firstValCol]
            }
            break
        }
    }
}

# --- Step 2: re-sort the table by "Casos totales" (column B) descending ---
$sortRange = $ws.Range("A$($firstDataRow):H$($lastDataRow)")
$sortKey = $ws.Range("B$firstDataRow")
$sortRange.Sort($sortKey, 2)

# --- Step 3: fix the Santa Lucia / Timor Oriental tie-break ---
$rowSantaLucia = -1
$rowTimorOriental = -1
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $existingName = $ws.Cells.Item($r, $nameCol).Value2
    if ($existingName -eq "Santa Lucia") { $rowSantaLucia = $r }
    if ($existingName -eq "Timor Oriental") { $rowTimorOriental = $r }
}

if ($rowSantaLucia -gt 0 -and $rowTimorOriental -gt 0 -and $rowSantaLucia -lt $rowTimorOriental) {
    for ($c = $nameCol; $c -le $lastValCol; $c++) {
        $valSanta = $ws.Cells.Item($rowSantaLucia, $c).Value2
        $valTimor = $ws.Cells.Item($rowTimorOriental, $c).Value2
        $ws.Cells.Item($rowSantaLucia, $c).Value = $valTimor
        $ws.Cells.Item($rowTimorOriental, $c).Value = $valSanta
    }
}

# --- Step 4: bump the "updated at" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 10:12"
